# Scheduled-runner style refresh of the Leve profit columns (H:N) across
# every job sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR). Each block below
# rewrites the currentAveragePrice* / LevePrice* / LeveProfit* cells for a
# single leve row with freshly pulled market-board figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 456.66666
$ws.Range("I2").Value = 150
$ws.Range("J2").Value = 610
$ws.Range("K2").Value = 150
$ws.Range("L2").Value = 610
$ws.Range("M2").Value = -37
$ws.Range("N2").Value = -836

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1815.3529
$ws.Range("I40").Value = 1840.88
$ws.Range("J40").Value = 1744.4445
$ws.Range("K40").Value = 1840.88
$ws.Range("L40").Value = 1744.4445
$ws.Range("M40").Value = -1665.88
$ws.Range("N40").Value = -2094.4445

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3964.6309
$ws.Range("I64").Value = 3426
$ws.Range("J64").Value = 5760.067
$ws.Range("K64").Value = 3426
$ws.Range("L64").Value = 5760.067
$ws.Range("M64").Value = -3178
$ws.Range("N64").Value = -6256.067

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3964.6309
$ws.Range("I67").Value = 3426
$ws.Range("J67").Value = 5760.067
$ws.Range("K67").Value = 3426
$ws.Range("L67").Value = 5760.067
$ws.Range("M67").Value = -2568
$ws.Range("N67").Value = -7476.067

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 950.6
$ws.Range("I98").Value = 866.0769
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 866.0769
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = 631.9231
$ws.Range("N98").Value = -4496

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 950.6
$ws.Range("I122").Value = 866.0769
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 2598.2307
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -148.2307000000001
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 38666.668
$ws.Range("J123").Value = 38666.668
$ws.Range("L123").Value = 38666.668
$ws.Range("N123").Value = -48466.668

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H126").Value = 12843.214
$ws.Range("J126").Value = 12843.214
$ws.Range("L126").Value = 12843.214
$ws.Range("N126").Value = -22723.214

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1002.39703
$ws.Range("I129").Value = 709.8
$ws.Range("J129").Value = 1025.619
$ws.Range("K129").Value = 2129.4
$ws.Range("L129").Value = 3076.857
$ws.Range("M129").Value = 2870.6
$ws.Range("N129").Value = -13076.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3847.32
$ws.Range("I61").Value = 2782.9048
$ws.Range("K61").Value = 2782.9048
$ws.Range("M61").Value = -2570.9048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 53711.6
$ws.Range("J133").Value = 53711.6
$ws.Range("L133").Value = 53711.6
$ws.Range("N133").Value = -58771.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3847.32
$ws.Range("I136").Value = 2782.9048
$ws.Range("K136").Value = 8348.714399999999
$ws.Range("M136").Value = -5798.714399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1546.3
$ws.Range("I99").Value = 968.46155
$ws.Range("J99").Value = 1988.1765
$ws.Range("K99").Value = 968.46155
$ws.Range("L99").Value = 1988.1765
$ws.Range("M99").Value = 529.53845
$ws.Range("N99").Value = -4984.1765

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2131.394
$ws.Range("I134").Value = 901.2727
$ws.Range("K134").Value = 2703.8181
$ws.Range("M134").Value = -168.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 20051
$ws.Range("J48").Value = 20051
$ws.Range("L48").Value = 20051
$ws.Range("N48").Value = -21003

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 32763.75
$ws.Range("J57").Value = 37000
$ws.Range("L57").Value = 37000
$ws.Range("N57").Value = -38120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2638.2104
$ws.Range("I58").Value = 1468.8572
$ws.Range("J58").Value = 3320.3333
$ws.Range("K58").Value = 1468.8572
$ws.Range("L58").Value = 3320.3333
$ws.Range("M58").Value = -1265.8572
$ws.Range("N58").Value = -3726.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2638.2104
$ws.Range("I136").Value = 1468.8572
$ws.Range("J136").Value = 3320.3333
$ws.Range("K136").Value = 4406.571599999999
$ws.Range("L136").Value = 9960.999899999999
$ws.Range("M136").Value = -1856.571599999999
$ws.Range("N136").Value = -15060.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1059.8966
$ws.Range("I5").Value = 631.1053000000001
$ws.Range("J5").Value = 1874.6
$ws.Range("K5").Value = 1893.3159
$ws.Range("L5").Value = 5623.799999999999
$ws.Range("M5").Value = -1781.3159
$ws.Range("N5").Value = -5847.799999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1059.8966
$ws.Range("I135").Value = 631.1053000000001
$ws.Range("J135").Value = 1874.6
$ws.Range("K135").Value = 5679.947700000001
$ws.Range("L135").Value = 16871.4
$ws.Range("M135").Value = -3144.947700000001
$ws.Range("N135").Value = -21941.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2200.9756
$ws.Range("I122").Value = 1911.0358
$ws.Range("J122").Value = 2825.4614
$ws.Range("K122").Value = 5733.107400000001
$ws.Range("L122").Value = 8476.3842
$ws.Range("M122").Value = -3283.107400000001
$ws.Range("N122").Value = -13376.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 202560.8
$ws.Range("I7").Value = 333933.34
$ws.Range("J7").Value = 5502
$ws.Range("K7").Value = 333933.34
$ws.Range("L7").Value = 5502
$ws.Range("M7").Value = -333821.34
$ws.Range("N7").Value = -5726

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H69").Value = 32000
$ws.Range("J69").Value = 32000
$ws.Range("L69").Value = 32000
$ws.Range("N69").Value = -33622

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H72").Value = 32000
$ws.Range("J72").Value = 32000
$ws.Range("L72").Value = 96000
$ws.Range("N72").Value = -104112

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 22680.6
$ws.Range("I93").Value = 28100.75
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 28100.75
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = -26852.75
$ws.Range("N93").Value = -3496

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 202560.8
$ws.Range("I126").Value = 333933.34
$ws.Range("J126").Value = 5502
$ws.Range("K126").Value = 1001800.02
$ws.Range("L126").Value = 16506
$ws.Range("M126").Value = -999330.02
$ws.Range("N126").Value = -21446

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 55000
$ws.Range("J58").Value = 55000
$ws.Range("L58").Value = 55000
$ws.Range("N58").Value = -55616

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 306.69232
$ws.Range("I107").Value = 275.18182
$ws.Range("J107").Value = 480
$ws.Range("K107").Value = 825.54546
$ws.Range("L107").Value = 1440
$ws.Range("M107").Value = 1094.45454
$ws.Range("N107").Value = -5280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 55424.684
$ws.Range("I126").Value = 69011.60000000001
$ws.Range("J126").Value = 4473.75
$ws.Range("K126").Value = 207034.8
$ws.Range("L126").Value = 13421.25
$ws.Range("M126").Value = -204564.8
$ws.Range("N126").Value = -18361.25
